$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.983.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.264.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.610.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.872"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.271.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.897.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0993"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.28%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0883"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.108"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +31.45%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.778.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.74%  "
